$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume/hour data as scraped at the new run time.
# Each target cell is stored as text (inline string) in the workbook, so we
# force the NumberFormat to Text before assigning, then restore the original
# style, to avoid Excel auto-coercing numeric-looking strings into numbers.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue ($ws.Range("D2")) '261.83'
Set-TextValue ($ws.Range("E2")) '1.73%'
Set-TextValue ($ws.Range("G2")) '1'
Set-TextValue ($ws.Range("D3")) '27.07'
Set-TextValue ($ws.Range("E3")) '0.23%'
Set-TextValue ($ws.Range("G3")) '1'
Set-TextValue ($ws.Range("D4")) '4.682'
Set-TextValue ($ws.Range("E4")) '0.87%'
Set-TextValue ($ws.Range("G4")) '1'
Set-TextValue ($ws.Range("D5")) '0.06081'
Set-TextValue ($ws.Range("E5")) '2.88%'
Set-TextValue ($ws.Range("G5")) '1'
Set-TextValue ($ws.Range("D6")) '6.712'
Set-TextValue ($ws.Range("E6")) '1.06%'
Set-TextValue ($ws.Range("G6")) '1'
Set-TextValue ($ws.Range("D7")) '0.8645'
Set-TextValue ($ws.Range("E7")) '1.14%'
Set-TextValue ($ws.Range("G7")) '1'
Set-TextValue ($ws.Range("D8")) '0.9207'
Set-TextValue ($ws.Range("E8")) '-2.04%'
Set-TextValue ($ws.Range("G8")) '1'
Set-TextValue ($ws.Range("D9")) '0.1410'
Set-TextValue ($ws.Range("E9")) '0.06%'
Set-TextValue ($ws.Range("G9")) '1'
Set-TextValue ($ws.Range("D10")) '0.04972'
Set-TextValue ($ws.Range("E10")) '10.76%'
Set-TextValue ($ws.Range("G10")) '1'
Set-TextValue ($ws.Range("D11")) '0.07129'
Set-TextValue ($ws.Range("E11")) '0.33%'
Set-TextValue ($ws.Range("G11")) '1'
Set-TextValue ($ws.Range("D12")) '0.03065'
Set-TextValue ($ws.Range("E12")) '-2.00%'
Set-TextValue ($ws.Range("G12")) '1'
Set-TextValue ($ws.Range("D13")) '0.09138'
Set-TextValue ($ws.Range("E13")) '-0.13%'
Set-TextValue ($ws.Range("G13")) '1'
Set-TextValue ($ws.Range("D14")) '0.001531'
Set-TextValue ($ws.Range("E14")) '0.57%'
Set-TextValue ($ws.Range("G14")) '1'
Set-TextValue ($ws.Range("D15")) '0.0006081'
Set-TextValue ($ws.Range("E15")) '0.21%'
Set-TextValue ($ws.Range("G15")) '1'
Set-TextValue ($ws.Range("D16")) '0.006193'
Set-TextValue ($ws.Range("E16")) '0.49%'
Set-TextValue ($ws.Range("G16")) '1'
Set-TextValue ($ws.Range("D17")) '3.495'
Set-TextValue ($ws.Range("E17")) '-0.39%'
Set-TextValue ($ws.Range("G17")) '1'
Set-TextValue ($ws.Range("D18")) '3.168'
Set-TextValue ($ws.Range("E18")) '-1.07%'
Set-TextValue ($ws.Range("G18")) '1'
Set-TextValue ($ws.Range("D19")) '2.198'
Set-TextValue ($ws.Range("E19")) '-0.26%'
Set-TextValue ($ws.Range("G19")) '1'
Set-TextValue ($ws.Range("D20")) '0.3128'
Set-TextValue ($ws.Range("G20")) '1'
Set-TextValue ($ws.Range("D21")) '0.1290'
Set-TextValue ($ws.Range("E21")) '-0.80%'
Set-TextValue ($ws.Range("G21")) '1'
Set-TextValue ($ws.Range("D22")) '4.096'
Set-TextValue ($ws.Range("E22")) '7.22%'
Set-TextValue ($ws.Range("G22")) '1'
Set-TextValue ($ws.Range("D23")) '0.04260'
Set-TextValue ($ws.Range("E23")) '0.05%'
Set-TextValue ($ws.Range("G23")) '1'
Set-TextValue ($ws.Range("E24")) '-0.11%'
Set-TextValue ($ws.Range("G24")) '1'
Set-TextValue ($ws.Range("D25")) '0.004056'
Set-TextValue ($ws.Range("E25")) '-5.65%'
Set-TextValue ($ws.Range("G25")) '1'
Set-TextValue ($ws.Range("E26")) '0.05%'
Set-TextValue ($ws.Range("G26")) '1'
Set-TextValue ($ws.Range("D27")) '0.0001564'
Set-TextValue ($ws.Range("E27")) '-19.24%'
Set-TextValue ($ws.Range("G27")) '1'
Set-TextValue ($ws.Range("G28")) '1'
Set-TextValue ($ws.Range("G29")) '1'
Set-TextValue ($ws.Range("G30")) '1'
Set-TextValue ($ws.Range("G31")) '1'
Set-TextValue ($ws.Range("G32")) '1'
Set-TextValue ($ws.Range("G33")) '1'
Set-TextValue ($ws.Range("G34")) '1'
Set-TextValue ($ws.Range("G35")) '1'
Set-TextValue ($ws.Range("G36")) '1'
Set-TextValue ($ws.Range("G37")) '1'
Set-TextValue ($ws.Range("G38")) '1'
Set-TextValue ($ws.Range("G39")) '1'
Set-TextValue ($ws.Range("D40")) '0.03883'
Set-TextValue ($ws.Range("E40")) '1.42%'
Set-TextValue ($ws.Range("G40")) '1'
Set-TextValue ($ws.Range("D41")) '0.1117'
Set-TextValue ($ws.Range("E41")) '1.19%'
Set-TextValue ($ws.Range("G41")) '1'
Set-TextValue ($ws.Range("D42")) '0.004151'
Set-TextValue ($ws.Range("E42")) '-33.98%'
Set-TextValue ($ws.Range("G42")) '1'
Set-TextValue ($ws.Range("D43")) '0.01518'
Set-TextValue ($ws.Range("E43")) '29.22%'
Set-TextValue ($ws.Range("G43")) '1'
Set-TextValue ($ws.Range("E44")) '0.06%'
Set-TextValue ($ws.Range("G44")) '1'
Set-TextValue ($ws.Range("D45")) '0.00005284'
Set-TextValue ($ws.Range("E45")) '-3.32%'
Set-TextValue ($ws.Range("G45")) '1'
Set-TextValue ($ws.Range("E46")) '0.06%'
Set-TextValue ($ws.Range("G46")) '1'
Set-TextValue ($ws.Range("E47")) '7.04%'
Set-TextValue ($ws.Range("G47")) '1'
Set-TextValue ($ws.Range("E48")) '-43.90%'
Set-TextValue ($ws.Range("G48")) '1'
Set-TextValue ($ws.Range("E49")) '0.06%'
Set-TextValue ($ws.Range("G49")) '1'
Set-TextValue ($ws.Range("E50")) '0.06%'
Set-TextValue ($ws.Range("G50")) '1'
Set-TextValue ($ws.Range("G51")) '1'
